# Insert a new row at position 24, shifting existing rows 24-27 down to 25-28.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(24).Insert()

# Populate the new row 24 with the new weekly record.
$ws.Cells.Item(24, 1).Value2 = 1
$ws.Cells.Item(24, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(24, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(24, 4).Value2 = 44769
$ws.Cells.Item(24, 5).Value2 = 15
$ws.Cells.Item(24, 6).Value2 = 100112013
$ws.Cells.Item(24, 7).Value2 = "Alcachofa"
$ws.Cells.Item(24, 8).Value2 = "Madrigal"
$ws.Cells.Item(24, 9).Value2 = "Primera"
$ws.Cells.Item(24, 10).Value2 = 200
$ws.Cells.Item(24, 11).Value2 = 17000
$ws.Cells.Item(24, 12).Value2 = 18000
$ws.Cells.Item(24, 13).Value2 = 17500
$ws.Cells.Item(24, 14).Value2 = "$/caja 40 unidades"
$ws.Cells.Item(24, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(24, 16).Value2 = 438
$ws.Cells.Item(24, 17).Value2 = 40
$ws.Cells.Item(24, 18).Value2 = "Hortaliza"

$wb.Save()
